$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 305; this shifts existing rows 305-381 down to 306-382
# and automatically extends the sheet dimension / UsedRange.
$ws.Rows.Item(305).Insert()

# Populate the newly inserted row 305 with the new record data.
$ws.Range("A305").Value = 5
$ws.Range("B305").Value = "Macroferia Regional de Talca"
$ws.Range("C305").Value = "Maule"
$ws.Range("D305").Value = 45135
$ws.Range("E305").Value = 7
$ws.Range("F305").Value = 100112021
$ws.Range("G305").Value = "Ají"
$ws.Range("H305").Value = "Inferno"
$ws.Range("I305").Value = "Primera"
$ws.Range("J305").Value = 150
$ws.Range("K305").Value = 12000
$ws.Range("L305").Value = 12000
$ws.Range("M305").Value = 12000
$ws.Range("N305").Value = "$/caja 10 kilos"
$ws.Range("O305").Value = "Región de Arica y Parinacota"
$ws.Range("P305").Value = 1200
$ws.Range("Q305").Value = 10
$ws.Range("R305").Value = "Hortaliza"
